# Apply the POAI 2025 staging updates:
#  - Introduce a new "Col19" column: the existing "Hoja" column (S) is
#    copied into a new column T (keeping its values/formatting), column S's
#    header is relabeled "Col19", and the old S data cells are cleared
#    (the "Hoja" values now live exclusively in column T).
#  - Update a handful of "Enlace" responsible-person values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Hoja" column from S to T (copies values + styles) -------
$ws.Range("S1:S15").Copy($ws.Range("T1:T15"))

# Relabel the old column S header, and clear the old S data cells.
$ws.Range("S1").Value = "Col19"
$ws.Range("S2:S15").ClearContents()

# --- Update responsible-person values ----------------------------------
$ws.Range("M6").Value = "ANDREA GONZALEZ"
$ws.Range("L8").Value = "SARA DIANA URBANO"
$ws.Range("L11").Value = "LUZ MIRYAN Y WILLIAN MOSQUERA"
